# Repull data, push all data, mean calculation
# Updates the dSF (F) column values for rows that were re-pulled/recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 0
$ws.Range("F3").Value  = 0
$ws.Range("F13").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("F17").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("F25").Value = -8
$ws.Range("F31").Value = -1
$ws.Range("F40").Value = -1
$ws.Range("F41").Value = -1
$ws.Range("F44").Value = 2
$ws.Range("F46").Value = 3
